$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 11 (currently "Flask-Uploads==0.2.1"), to add Flask-SQLAlchemy after Flask
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Flask-SQLAlchemy==2.3.2"

# After the previous insert, "six==1.12.0" (was row 28) is now row 29.
# Insert a new row before row 30 (currently "Werkzeug==0.14.1"), to add SQLAlchemy after six
$ws.Rows.Item(30).Insert()
$ws.Range("A30").Value = "SQLAlchemy==1.3.0"
